$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capitalize the "sd" header abbreviation used in C3 / E3 -> "SD"
$ws.Range("C3").Value = "SD"
$ws.Range("E3").Value = "SD"

# Rename the title cell (A1): "Mini-gland structure measurement" -> "Parotid gland measurement"
$ws.Range("A1").Value = "Parotid gland measurement"

# Adjust column widths (A, C, F) to the new layout
$ws.Columns.Item(1).ColumnWidth = 23
$ws.Columns.Item(3).ColumnWidth = 4.333333333333334
$ws.Columns.Item(6).ColumnWidth = 31.333333333333332

# Move the active selection from F21 down to the full next row (A22:XFD22)
$ws.Rows.Item(22).Select()
